$wb = $excel.ActiveWorkbook

# --- Sheet "具有相當價值之財產" -> "保險" ---------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Name = "保險"

# The sheet used to contain a duplicated header-label row (row 2, index 125)
# immediately above the real data rows. Remove that duplicate row so the
# 4 real insurance records (indices 126-129) move up to rows 2-5.
$ws.Rows(2).Delete()

# Row 1 (previously the column-header labels 保險公司/保險名稱/要保人/備註)
# now mirrors the first insurance record instead of header text.
$ws.Range("B1").Value = "國泰人壽"
$ws.Range("C1").Value = "保本111終身"
$ws.Range("D1").Value = "陳淑慧"
$ws.Range("E1").Value = ""
